# Add frog to master data.
# Appends a new "chr_frog" row to the "Character" sheet, mirroring the
# existing rows (id=3, assetName="chr_frog", moveSpeed=400, weight=100,
# jumpHeight=200, aerialJumpCount=1, power=2, hp=3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Character")

# Copy the formatting of the previous data row (row 6) down into the new
# row (row 7) so the new row matches the existing bordered style.
$ws.Range("A6:H6").Copy()
$ws.Range("A7:H7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = "chr_frog"
$ws.Cells.Item(7, 3).Value = 400
$ws.Cells.Item(7, 4).Value = 100
$ws.Cells.Item(7, 5).Value = 200
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2
$ws.Cells.Item(7, 8).Value = 3

$ws.Range("F21").Select()
